$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nome = "DADOS FICTICIOS APENAS PARA FINS DE SIMULAÇÃO"
$vinculo = "VINCULO TESTE"

$cpfs = @(
    20412634524,
    10641944494,
    11791847506,
    16785564396,
    16862797558,
    10224458768,
    21447448571,
    16553730149,
    11742077000,
    12017882684,
    16863150477,
    16430525776,
    12798121452,
    16398233752,
    12580267486,
    16146164616,
    20453212950,
    16361202381,
    20767551642,
    13378191421,
    16400283375,
    12493595779,
    12944297653,
    16140231826
)

$startRow = 11
for ($i = 0; $i -lt $cpfs.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $nome
    $ws.Cells.Item($row, 2).Value = 99999999999
    $ws.Cells.Item($row, 3).Value = $cpfs[$i]
    $ws.Cells.Item($row, 4).Value = $vinculo
    $ws.Cells.Item($row, 5).Value = 999999999.99
}

$ws.Range("A32").Select()
